{"js": "// Replace the constellation name \"Perseus\" with \"Cygnus constellation\" in the\n// first introductory paragraph (\"You are participating in a global\n// campaign...\"), matching the author's edit which also collapsed that\n// paragraph's many small runs into a single run of plain text.\n\nconst body = context.document.body;\n\n// Locate the paragraph that mentions \"Perseus\" (the constellation name being\n// swapped out) instead of relying on a hard-coded paragraph index.\nconst searchResults = body.search(\"Perseus\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the target paragraph (text 'Perseus' not found).\");\n}\n\nconst paragraph = searchResults.items[0].paragraphs.getFirst();\nparagraph.load(\"text\");\nawait context.sync();\n\nconst newText = paragraph.text.replace(\"Perseus\", \"Cygnus constellation\");\n\n// Clear the paragraph's existing (many-run) content, then insert the full,\n// updated sentence as a single run.\nparagraph.clear();\nawait context.sync();\n\nparagraph.insertText(newText, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Replace the constellation name \"Perseus\" with \"Cygnus constellation\" in the\n# first introductory paragraph (\"You are participating in a global\n# campaign...\"), matching the author's edit which also collapsed that\n# paragraph's many small runs into a single run of plain text.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that mentions \"Perseus\" (the constellation name being\n# swapped out) instead of relying on a hard-coded paragraph index.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Perseus\"\n$find.Forward = $true\n$find.Wrap = 0\n$found = $find.Execute()\nif (-not $found) {\n  throw \"Could not find the target paragraph (text 'Perseus' not found).\"\n}\n\n# Expand the found hit out to its whole paragraph, excluding the trailing\n# paragraph mark so the replacement stays inside this paragraph.\n$rng.Expand(4) | Out-Null   # wdParagraph\n$rng.End = $rng.End - 1\n\n$newText = $rng.Text.Replace(\"Perseus\", \"Cygnus constellation\")\n\n# Delete the old (many-run) paragraph content and insert the updated\n# sentence back in as a single plain run.\n$rng.Delete()\n$rng.InsertAfter($newText)\n"}
